# Update the embedded build timestamp throughout the workbook:
#   "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"
# This timestamp appears in the "About" sheet (version line + citation line)
# and in the "build_version" column (S2:S10) of the
# "Boundaries and methane sources" sheet.

$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find($oldTimestamp)
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        while ($true) {
            $newText = $found.Text.Replace($oldTimestamp, $newTimestamp)
            $found.Value = $newText
            $found = $used.FindNext($found)
            if ($found -eq $null) {
                break
            }
            if ($found.Address() -eq $firstAddress) {
                break
            }
        }
    }
}
